$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repull/recalculation of data
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 6
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = -2
